$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emp_data")

$ws.Range("A6").Value = "RE1474"
$ws.Range("B6").Value = "swamiraj"
$ws.Range("C6").Value = "swami@reposenergy.com"
$ws.Range("D6").Value = "M"

$ws.Range("A7").Value = "RE6472"
$ws.Range("B7").Value = "akshay"
$ws.Range("C7").Value = "aksh@gmail.com"
$ws.Range("D7").Value = "M"
